$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "You need the asa.sql.admin password. My suggestion was P@ssw0rd1234. If someone used another one and forgot it, then reset the password in the Synapse Workspace in the portal"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "You need the asa.sql.admin password. My suggestion was P@ssw0rd. If someone used another one and forgot it, then reset the password in the Synapse Workspace in the portal"
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "Don’t forget to pause the sql pool at the end"
$find2.Replacement.ClearFormatting()
$find2.Replacement.Text = "Don’t forget to pause the sql pool at the end"
$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null
